$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a value into $cellRef as a genuine text cell, without leaving
# any NumberFormat/style change behind on the destination cell. A
# scratch cell (Z1, well outside the A1:E51 table) is forced to Text,
# given the value, then copied with Paste Special (values only) onto
# the real destination -- this carries over the TEXT storage type
# without carrying over Z1's own formatting.
function Set-TextValue($cellRef, $val) {
    $scratch = $ws.Range("Z1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# --- Column D ("Price") + column E ("Volume(1h)") updates ---
Set-TextValue "D2" "34.090.55"
$ws.Range("E2").Value = "  -1.66%  "
Set-TextValue "D3" "1.787.34"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "222.83"
$ws.Range("E5").Value = "  -1.14%  "
Set-TextValue "D6" "0.550"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue "D8" "32.32"
$ws.Range("E8").Value = "  -1.19%  "
Set-TextValue "D9" "0.285"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  +0.21%  "
Set-TextValue "D12" "2.044.00"
Set-TextValue "D13" "1.806.55"
$ws.Range("E13").Value = "  -0.14%  "
Set-TextValue "D14" "10.93"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("E15").Value = "  -3.09%  "
Set-TextValue "D16" "34.067.93"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("E17").Value = "  -3.99%  "
Set-TextValue "D18" "67.94"
$ws.Range("E18").Value = "  -2.38%  "
Set-TextValue "D19" "243.62"
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("E24").Value = "  -2.49%  "
Set-TextValue "D25" "158.61"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("E33").Value = "  -4.26%  "
$ws.Range("E34").Value = "  -4.90%  "
Set-TextValue "D35" "1.383.26"
$ws.Range("E35").Value = "  -3.88%  "
Set-TextValue "D36" "0.647"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  -4.33%  "
Set-TextValue "D39" "79.56"
$ws.Range("E39").Value = "  -6.48%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -4.52%  "
Set-TextValue "D42" "2.70"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("E43").Value = "  +0.17%  "
Set-TextValue "D44" "0.0₆0137"
$ws.Range("E44").Value = "  +6.97%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  -0.71%  "
Set-TextValue "D47" "107.35"
$ws.Range("E47").Value = "  +0.93%  "
Set-TextValue "D48" "5.84"
$ws.Range("E48").Value = "  -3.19%  "
Set-TextValue "D49" "1.943.91"
$ws.Range("E49").Value = "  -0.85%  "

# --- Rows 50/51: PaxDollar <-> InjectiveProtocol swap positions ---
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D50" "11.99"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  -0.03%  "

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$ws.Range("Z1").Clear()
